# This script updates the "Pedido" order quantities (column L, "Unidades Pedido")
# and "Diferencia Stock" (column M) for a set of article rows so that the
# "Diferencia Stock" becomes 0 for those rows (i.e. the extra stock-difference
# units are removed from the order), and refreshes the summary metrics at the
# bottom of the sheet (resumen_pedido) accordingly: Total_Unidades (C139) and
# Total_Ajuste_Stock (C150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (new L value, new M value)
$updates = @{
    13  = @(3, 0)
    14  = @(4, 0)
    16  = @(3, 0)
    20  = @(4, 0)
    21  = @(7, 0)
    24  = @(9, 0)
    25  = @(5, 0)
    26  = @(9, 0)
    27  = @(4, 0)
    30  = @(13, 0)
    31  = @(6, 0)
    47  = @(1, 0)
    56  = @(7, 0)
    65  = @(3, 0)
    72  = @(1, 0)
    74  = @(1, 0)
    79  = @(4, 0)
    87  = @(6, 0)
    96  = @(7, 0)
    100 = @(1, 0)
    105 = @(1, 0)
    109 = @(4, 0)
    110 = @(7, 0)
    114 = @(1, 0)
    123 = @(1, 0)
    128 = @(2, 0)
    131 = @(1, 0)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("L$row").Value2 = $vals[0]
    $ws.Range("M$row").Value2 = $vals[1]
}

# Update the summary metrics section ("METRICAS DE RESUMEN")
$ws.Range("C139").Value2 = 334   # Total_Unidades
$ws.Range("C150").Value2 = 0     # Total_Ajuste_Stock
